$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fecha" (date) column B is being dropped; "Descuento" and "Porcentaje"
# (previously in C and D) move one column to the left, into B and C.

# Row 1 - headers
$ws.Range("B1").Value = "Descuento"
$ws.Range("C1").Value = "Porcentaje"

# Row 2 - Falta
$ws.Range("B2").Value = 84.5
$ws.Range("C2").Value = 0.9

# Row 3 - Comida
$ws.Range("B3").Value = 95.62
$ws.Range("C3").Value = 0.5

# Row 4 - Gasolina
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 0.5

# Column D held "Porcentaje" before the shift and is no longer used; clear it.
$ws.Range("D1:D4").Clear()

# Column B held dates (with a date number format); strip that formatting now
# that it holds the plain "Descuento" numbers/header instead.
$ws.Range("B2:B4").ClearFormats()

# Match the author's final selection in the saved file.
$ws.Range("D3").Select()
